# ------------------------------------------------------------------
# Applies the "additional scraping" edit:
#   1. Inserts a new "Player Info" worksheet as the first sheet,
#      containing ID / NAME / BATTING_HAND / BOWL_STYLE for the player.
#   2. Renames MATCH_CARD_LINK -> MATCH_CODE on both the "ODI Batting"
#      and "ODI Bowling" sheets, replacing the full scorecard URL
#      values with just the numeric match code (kept as text).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$batting = $wb.Worksheets.Item("ODI Batting")
$bowling = $wb.Worksheets.Item("ODI Bowling")

# --- 1. Create the new "Player Info" sheet in front of "ODI Batting" ---
$playerInfo = $wb.Worksheets.Add($batting)
$playerInfo.Name = "Player Info"

# Re-fetch sheet references: inserting a sheet shifts worksheet
# positions, and stale references end up acting on the wrong sheet.
$batting = $wb.Worksheets.Item("ODI Batting")
$bowling = $wb.Worksheets.Item("ODI Bowling")
$playerInfo = $wb.Worksheets.Item("Player Info")

# --- Populate "Player Info" ---
$playerInfo.Activate()

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Re-use the same bold/centered/bordered header formatting as the
# other sheets' header row (style index 1 in the original workbook).
$batting.Activate()
$batting.Range("A1").Copy()
$playerInfo.Activate()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4488"
$playerInfo.Range("B2").Value = "Iftikhar Ahmed"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# --- 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (column D) ---
$batting.Activate()

$battingCodes = @("3859", "3861", "4375", "4376", "4432", "4433", "4434", "4564", "4565", "4567")

$batting.Range("D1").Value = "MATCH_CODE"

for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $row = $i + 2
    $cell = $batting.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $battingCodes[$i]
}

# --- 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE (column B) ---
$bowling.Activate()

$bowlingCodes = @("3859", "3861", "4375", "4376", "4433", "4434", "4564", "4565", "4567")

$bowling.Range("B1").Value = "MATCH_CODE"

for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $row = $i + 2
    $cell = $bowling.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $bowlingCodes[$i]
}

# Restore the first sheet as the active tab, same as the original file.
$playerInfo.Activate()
$playerInfo.Range("A1").Select() | Out-Null

$wb.Save()
